# "Nuevo formato 15 jun 2021"
# Update the statistics tables (Blancos, Aprobados, Por_Apro, Promedio) on
# the "Estadisticos 1P", "Estadisticos 2P" and "Estadisticos Final" sheets.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "Estadisticos 1P" = @(
        @{ Row = 2; D = 1; F = 20; G = 95.23999999999999; H = 8.699999999999999 },
        @{ Row = 3; D = 1; F = 20; G = 95.23999999999999; H = 8.6 },
        @{ Row = 4; D = 4; F = 35; G = 89.73999999999999; H = 7.6 },
        @{ Row = 5; D = 2; F = 19; G = 90.48;              H = 8.699999999999999 }
    )
    "Estadisticos 2P" = @(
        @{ Row = 2; D = 1; F = 20; G = 95.23999999999999; H = 8.699999999999999 },
        @{ Row = 3; D = 1; F = 20; G = 95.23999999999999; H = 8.6 },
        @{ Row = 4; D = 4; F = 35; G = 89.73999999999999; H = 7.6 },
        @{ Row = 5; D = 2; F = 19; G = 90.48;              H = 8.699999999999999 }
    )
    "Estadisticos Final" = @(
        @{ Row = 2; D = 1; F = 20; G = 95.23999999999999; H = 8 },
        @{ Row = 3; D = 1; F = 20; G = 95.23999999999999; H = 8.199999999999999 },
        @{ Row = 4; D = 4; F = 35; G = 89.73999999999999; H = 7.5 },
        @{ Row = 5; D = 2; F = 19; G = 90.48;              H = 8.9 }
    )
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($update in $sheetUpdates[$sheetName]) {
        $r = $update.Row
        $ws.Cells.Item($r, 4).Value = $update.D
        $ws.Cells.Item($r, 6).Value = $update.F
        $ws.Cells.Item($r, 7).Value = $update.G
        $ws.Cells.Item($r, 8).Value = $update.H
    }
}
